$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 28 (which holds problem 128,
# "最长连续序列"). This pushes that row down to row 30 and opens up rows
# 28-29 for the two newly solved problems (27 & 28).
$ws.Rows("28:29").Insert()

# Carry the per-column formatting (borders / number formats / fonts) from
# the row above down into the two new rows, without touching the rest of
# the (very wide) row so the file doesn't balloon with style overrides.
$ws.Range("A27:I27").Copy()
$ws.Range("A28:I29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 28: problem 27 - "移除元素" (remove-element)
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "简单"
$ws.Cells.Item(28, 3).Value = "移除元素"
$ws.Cells.Item(28, 4).Value = "remove-element"
$ws.Cells.Item(28, 5).Value = "191118-1.cpp"
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 8.4
$ws.Cells.Item(28, 8).Value = "100.00%"
$ws.Cells.Item(28, 9).Value = 36942276

# Row 29: problem 28 - "实现 strStr()" (implement-strstr)
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "简单"
$ws.Cells.Item(29, 3).Value = "实现 strStr()"
$ws.Cells.Item(29, 4).Value = "implement-strstr"
$ws.Cells.Item(29, 5).Value = "191118-1.cpp"
$ws.Cells.Item(29, 6).Value = 8
$ws.Cells.Item(29, 7).Value = 8.8
$ws.Cells.Item(29, 8).Value = "65.58%"
$ws.Cells.Item(29, 9).Value = 36942358
